# Trading update: 2026-02-17 15:48:28
# Appends the newly-closed/open "MarketMaking" trade (Trade # 71) as row 72
# to both the "All Trades" log and the per-strategy "MarketMaking" log.

$wb = $excel.ActiveWorkbook

$tradeNum      = 71
$tradeDate     = "2026-02-17"
$tradeTime     = "15:48:23"
$strategy      = "MarketMaking"
$side          = "UP"
$entryPrice    = 0.72
$exitPrice     = ""
$status        = "OPEN"
$pnlPct        = 0
$pnlDollar     = 0
$capitalAfter  = 100.2103383789229
$entrySlippage = 0
$exitSlippage  = 0
$confidence    = 0.6
$entryReason   = "Normal spread capture: 19600 bps"
$exitReason    = ""
$durationMin   = 0

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = $ws.UsedRange.Rows.Count + 1

    $ws.Cells.Item($row, 1).Value  = $tradeNum
    # Prefix date/time-shaped text with a leading apostrophe so the engine
    # keeps it as literal text instead of inferring a date/time serial.
    $ws.Cells.Item($row, 2).Value  = "'" + $tradeDate
    $ws.Cells.Item($row, 3).Value  = $tradeTime
    $ws.Cells.Item($row, 4).Value  = $strategy
    $ws.Cells.Item($row, 5).Value  = $side
    $ws.Cells.Item($row, 6).Value  = $entryPrice
    $ws.Cells.Item($row, 7).Value  = $exitPrice
    $ws.Cells.Item($row, 8).Value  = $status
    $ws.Cells.Item($row, 9).Value  = $pnlPct
    $ws.Cells.Item($row, 10).Value = $pnlDollar
    $ws.Cells.Item($row, 11).Value = $capitalAfter
    $ws.Cells.Item($row, 12).Value = $entrySlippage
    $ws.Cells.Item($row, 13).Value = $exitSlippage
    $ws.Cells.Item($row, 14).Value = $confidence
    $ws.Cells.Item($row, 15).Value = $entryReason
    $ws.Cells.Item($row, 16).Value = $exitReason
    $ws.Cells.Item($row, 17).Value = $durationMin
}
